$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell formatting (numFmt) from C5 onto C6:C8 so the new
# date cells reuse the existing short-date style instead of creating a
# brand new one.
$ws.Range("C5").Copy()
$ws.Range("C6:C8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 6 - D12 - 16 Jan 2020
$ws.Range("B6").Value = "D12"
$ws.Range("C6").Value = 43846
$ws.Range("D6").Value = "Reached lesson 5.4"

# Row 7 - D13 - 17 Jan 2020
$ws.Range("B7").Value = "D13"
$ws.Range("C7").Value = 43847
$ws.Range("D7").Value = "Reached lesson 5.5"

# Row 8 - D14 - 18 Jan 2020
$ws.Range("B8").Value = "D14"
$ws.Range("C8").Value = 43848
$ws.Range("D8").Value = "Completed Lesson 5.8, working on 5.9"

# Move the active selection to D8, matching where the author last edited.
$ws.Range("D8").Select()
